# Update "Nädal 4" worksheet (active sheet) per commit:
# "Update MeasuresPage, Create, Delete, details, Edit, Index pages"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 4")

# Row 15: fill in stop time, interruption time, delta time, comment (H15) and C-column mark (J15)
$ws.Range("D15").Value = 0.71388888888888891
$ws.Range("E15").Value = 60
$ws.Range("F15").Value = 240
$ws.Range("H15").Value = "p. 3"
$ws.Range("J15").Value = "x"

# Row 16: fill in date, start time, activity
$ws.Range("B16").Value = 43885
$ws.Range("C16").Value = 0.97222222222222221
$ws.Range("G16").Value = "Kodutöö 4"

# Update selected cell to C17 as per diff
$ws.Range("C17").Select()

$wb.Save()
